$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = "❌ EXPIRED 3757 days ago"
$ws.Range("E6").Value = "⚠️ Expires in 22 days"
$ws.Range("E7").Value = "⚠️ Expires in 22 days"
$ws.Range("E8").Value = "⚠️ Expires in 22 days"
